# This script adds a new "time_taken" column (column F) to the active
# worksheet. It writes the header "time_taken" into F1 (matching the style
# of the existing header row, i.e. the same style as A1:E1) and fills
# F2:F83 with per-row timestamp strings recorded when each panel row was
# captured.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Auto-generated list of time_taken values for F2:F83
$timeTakenValues = @(
    "2021-10-05 10:50:39.945343"
    "2021-10-05 10:50:39.945354"
    "2021-10-05 10:50:39.945357"
    "2021-10-05 10:50:39.945360"
    "2021-10-05 10:50:39.945363"
    "2021-10-05 10:50:39.945366"
    "2021-10-05 10:50:39.945368"
    "2021-10-05 10:50:39.945371"
    "2021-10-05 10:50:39.945373"
    "2021-10-05 10:50:39.945376"
    "2021-10-05 10:50:39.945378"
    "2021-10-05 10:50:39.945381"
    "2021-10-05 10:50:39.945383"
    "2021-10-05 10:50:39.945386"
    "2021-10-05 10:50:39.945388"
    "2021-10-05 10:50:39.945391"
    "2021-10-05 10:50:39.945394"
    "2021-10-05 10:50:39.945396"
    "2021-10-05 10:50:39.945399"
    "2021-10-05 10:50:39.945401"
    "2021-10-05 10:50:39.945404"
    "2021-10-05 10:50:39.945406"
    "2021-10-05 10:50:39.945409"
    "2021-10-05 10:50:39.945411"
    "2021-10-05 10:50:39.945414"
    "2021-10-05 10:50:39.945416"
    "2021-10-05 10:50:39.945419"
    "2021-10-05 10:50:39.945421"
    "2021-10-05 10:50:39.945424"
    "2021-10-05 10:50:39.945426"
    "2021-10-05 10:50:39.945429"
    "2021-10-05 10:50:39.945431"
    "2021-10-05 10:50:39.945434"
    "2021-10-05 10:50:39.945437"
    "2021-10-05 10:50:39.945440"
    "2021-10-05 10:50:39.945442"
    "2021-10-05 10:50:39.945445"
    "2021-10-05 10:50:39.945447"
    "2021-10-05 10:50:39.945450"
    "2021-10-05 10:50:39.945452"
    "2021-10-05 10:50:39.945455"
    "2021-10-05 10:50:39.945458"
    "2021-10-05 10:50:39.945460"
    "2021-10-05 10:50:39.945463"
    "2021-10-05 10:50:39.945465"
    "2021-10-05 10:50:39.945467"
    "2021-10-05 10:50:39.945470"
    "2021-10-05 10:50:39.945472"
    "2021-10-05 10:50:39.945475"
    "2021-10-05 10:50:39.945477"
    "2021-10-05 10:50:39.945480"
    "2021-10-05 10:50:39.945482"
    "2021-10-05 10:50:39.945485"
    "2021-10-05 10:50:39.945488"
    "2021-10-05 10:50:39.945490"
    "2021-10-05 10:50:39.945493"
    "2021-10-05 10:50:39.945495"
    "2021-10-05 10:50:39.945498"
    "2021-10-05 10:50:39.945500"
    "2021-10-05 10:50:39.945502"
    "2021-10-05 10:50:39.945505"
    "2021-10-05 10:50:39.945507"
    "2021-10-05 10:50:39.945510"
    "2021-10-05 10:50:39.945512"
    "2021-10-05 10:50:39.945516"
    "2021-10-05 10:50:39.945519"
    "2021-10-05 10:50:39.945521"
    "2021-10-05 10:50:39.945524"
    "2021-10-05 10:50:39.945526"
    "2021-10-05 10:50:39.945529"
    "2021-10-05 10:50:39.945531"
    "2021-10-05 10:50:39.945534"
    "2021-10-05 10:50:39.945536"
    "2021-10-05 10:50:39.945539"
    "2021-10-05 10:50:39.945541"
    "2021-10-05 10:50:39.945544"
    "2021-10-05 10:50:39.945549"
    "2021-10-05 10:50:39.945552"
    "2021-10-05 10:50:39.945554"
    "2021-10-05 10:50:39.945557"
    "2021-10-05 10:50:39.945559"
    "2021-10-05 10:50:39.945562"
)

# Header cell F1 - copy the formatting from the last existing header cell
# (E1) so that F1 matches the other header cells (bold, centered,
# bordered), then overwrite the copied value with the new header text.
$ws.Range("E1").Copy($ws.Range("F1"))
$ws.Range("F1").Value = "time_taken"

# Fill F2:F83 with the recorded time_taken values, one per data row.
$startRow = 2
for ($i = 0; $i -lt $timeTakenValues.Count; $i++) {
    $row = $startRow + $i
    $ws.Cells.Item($row, 6).Value = $timeTakenValues[$i]
}
